# Applies the FIN13Inicial.xlsx edit described in the commit "adding more
# space in semillerostable": widens the cronograma month headers from
# F/M/A/J (4 months) to A/S/O/N/D (5 months), and swaps out some of the
# placeholder/test data for real-looking sample data.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Hoja1. Actividades" ----
$ws1 = $wb.Worksheets.Item("Hoja1. Actividades")

# Group/seedbed name
$ws1.Range("E4").Value = "Semillero de Investigación en Inteligencia Computacional"

# Director name
$ws1.Range("E5").Value = "José Alejandro Cortés Taborda"

# Semester
$ws1.Range("E6").Value = "2019-2"

# Cronograma month headers: F, M, A, J -> A, S, O, N, D
$ws1.Range("D9").Value = "A"
$ws1.Range("E9").Value = "S"
$ws1.Range("F9").Value = "O"
$ws1.Range("G9").Value = "N"
$ws1.Range("H9").Value = "D"

# Row 10 (Actividad 1 schedule) - remove the "X" mark in column F
$ws1.Range("F10").ClearContents()

# Product name
$ws1.Range("J10").Value = " pActividad1"

# ---- Sheet "Hoja2. Integrantes" ----
$ws2 = $wb.Worksheets.Item("Hoja2. Integrantes")

$ws2.Range("A6").Value = "Usuario nuevo Prueba"
$ws2.Range("C6").Value = 1002
$ws2.Range("D6").Value = 1002
$ws2.Range("E6").Value = "Usuario_prueba@elpoli.edu.co"
